$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 320.25
$ws.Cells.Item(2, 9).Value = 640.2
$ws.Cells.Item(2, 10).Value = 91.71429000000001
$ws.Cells.Item(2, 11).Value = 640.2
$ws.Cells.Item(2, 12).Value = 91.71429000000001
$ws.Cells.Item(2, 13).Value = -527.2
$ws.Cells.Item(2, 14).Value = -317.71429

$ws.Cells.Item(21, 8).Value = 3800
$ws.Cells.Item(21, 10).Value = 3800
$ws.Cells.Item(21, 12).Value = 3800
$ws.Cells.Item(21, 14).Value = -4736

$ws.Cells.Item(23, 8).Value = 3800
$ws.Cells.Item(23, 10).Value = 3800
$ws.Cells.Item(23, 12).Value = 3800
$ws.Cells.Item(23, 14).Value = -4268

$ws.Cells.Item(33, 8).Value = 84.1579
$ws.Cells.Item(33, 9).Value = 63.384617
$ws.Cells.Item(33, 10).Value = 129.16667
$ws.Cells.Item(33, 11).Value = 63.384617
$ws.Cells.Item(33, 12).Value = 129.16667
$ws.Cells.Item(33, 13).Value = 165.615383
$ws.Cells.Item(33, 14).Value = -587.1666700000001

$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 13).Value = $null

$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 13).Value = $null

$ws.Cells.Item(70, 8).Value = 4866
$ws.Cells.Item(70, 10).Value = 6998
$ws.Cells.Item(70, 12).Value = 20994
$ws.Cells.Item(70, 14).Value = -21534

$ws.Cells.Item(73, 8).Value = 4866
$ws.Cells.Item(73, 10).Value = 6998
$ws.Cells.Item(73, 12).Value = 20994
$ws.Cells.Item(73, 14).Value = -22866

$ws.Cells.Item(80, 8).Value = 1655
$ws.Cells.Item(80, 10).Value = 1630
$ws.Cells.Item(80, 12).Value = 4890
$ws.Cells.Item(80, 14).Value = -6886

$ws.Cells.Item(83, 8).Value = 1655
$ws.Cells.Item(83, 10).Value = 1630
$ws.Cells.Item(83, 12).Value = 14670
$ws.Cells.Item(83, 14).Value = -24654

$ws.Cells.Item(88, 8).Value = 5415.6665
$ws.Cells.Item(88, 10).Value = 5415.6665
$ws.Cells.Item(88, 12).Value = 5415.6665
$ws.Cells.Item(88, 14).Value = -6227.6665

$ws.Cells.Item(91, 8).Value = 5415.6665
$ws.Cells.Item(91, 10).Value = 5415.6665
$ws.Cells.Item(91, 12).Value = 5415.6665
$ws.Cells.Item(91, 14).Value = -8223.666499999999

$ws.Cells.Item(95, 8).Value = 60000
$ws.Cells.Item(95, 10).Value = 60000
$ws.Cells.Item(95, 12).Value = 60000
$ws.Cells.Item(95, 14).Value = -65492

$ws.Cells.Item(96, 8).Value = 1835.625
$ws.Cells.Item(96, 10).Value = 3016.25
$ws.Cells.Item(96, 12).Value = 9048.75
$ws.Cells.Item(96, 14).Value = -11794.75

$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 13).Value = $null

$ws.Cells.Item(106, 8).Value = 2249.5
$ws.Cells.Item(106, 9).Value = 2249.5
$ws.Cells.Item(106, 11).Value = 2249.5
$ws.Cells.Item(106, 13).Value = -1618.5

$ws.Cells.Item(132, 8).Value = 3000
$ws.Cells.Item(132, 9).Value = 3000
$ws.Cells.Item(132, 11).Value = 9000
$ws.Cells.Item(132, 13).Value = -6470

$ws.Cells.Item(137, 8).Value = 4950
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).Value = $null

$ws.Cells.Item(138, 8).Value = 3412.9473
$ws.Cells.Item(138, 9).Value = 3309.4
$ws.Cells.Item(138, 11).Value = 9928.200000000001
$ws.Cells.Item(138, 13).Value = -4788.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2750
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 2750
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 2750
$ws.Cells.Item(2, 13).Value = $null
$ws.Cells.Item(2, 14).Value = -2976

$ws.Cells.Item(61, 8).Value = 15000
$ws.Cells.Item(61, 9).Value = 15000
$ws.Cells.Item(61, 11).Value = 15000
$ws.Cells.Item(61, 13).Value = -14788

$ws.Cells.Item(63, 8).Value = 4504.6665
$ws.Cells.Item(63, 9).Value = 4504.6665
$ws.Cells.Item(63, 11).Value = 4504.6665
$ws.Cells.Item(63, 13).Value = -3818.6665

$ws.Cells.Item(66, 8).Value = 4504.6665
$ws.Cells.Item(66, 9).Value = 4504.6665
$ws.Cells.Item(66, 11).Value = 22523.3325
$ws.Cells.Item(66, 13).Value = -19091.3325

$ws.Cells.Item(94, 8).Value = 29999.5
$ws.Cells.Item(94, 10).Value = 29999.5
$ws.Cells.Item(94, 12).Value = 29999.5
$ws.Cells.Item(94, 14).Value = -31801.5

$ws.Cells.Item(116, 8).Value = 2750
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 2750
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 2750
$ws.Cells.Item(116, 13).Value = $null
$ws.Cells.Item(116, 14).Value = -7338

$ws.Cells.Item(122, 8).Value = 1084.3334
$ws.Cells.Item(122, 10).Value = 1009.5
$ws.Cells.Item(122, 12).Value = 3028.5
$ws.Cells.Item(122, 14).Value = -7928.5

$ws.Cells.Item(132, 8).Value = 2560
$ws.Cells.Item(132, 9).Value = 2560
$ws.Cells.Item(132, 11).Value = 7680
$ws.Cells.Item(132, 13).Value = -5150

$ws.Cells.Item(136, 8).Value = 15000
$ws.Cells.Item(136, 9).Value = 15000
$ws.Cells.Item(136, 11).Value = 45000
$ws.Cells.Item(136, 13).Value = -42450

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2750
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 2750
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 2750
$ws.Cells.Item(3, 13).Value = $null
$ws.Cells.Item(3, 14).Value = -2978

$ws.Cells.Item(54, 8).Value = 3029.6667
$ws.Cells.Item(54, 9).Value = 3029.6667
$ws.Cells.Item(54, 11).Value = 3029.6667
$ws.Cells.Item(54, 13).Value = -2545.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 1500
$ws.Cells.Item(12, 9).Value = 1000
$ws.Cells.Item(12, 10).Value = 2000
$ws.Cells.Item(12, 11).Value = 1000
$ws.Cells.Item(12, 12).Value = 2000
$ws.Cells.Item(12, 13).Value = -830
$ws.Cells.Item(12, 14).Value = -2340

$ws.Cells.Item(22, 8).Value = 273.66666
$ws.Cells.Item(22, 10).Value = 360.5
$ws.Cells.Item(22, 12).Value = 360.5
$ws.Cells.Item(22, 14).Value = -1060.5

$ws.Cells.Item(31, 8).Value = 12500
$ws.Cells.Item(31, 10).Value = 12500
$ws.Cells.Item(31, 12).Value = 12500
$ws.Cells.Item(31, 14).Value = -13090

$ws.Cells.Item(34, 8).Value = 12500
$ws.Cells.Item(34, 10).Value = 12500
$ws.Cells.Item(34, 12).Value = 12500
$ws.Cells.Item(34, 14).Value = -12904

$ws.Cells.Item(35, 8).Value = 1716.1666
$ws.Cells.Item(35, 9).Value = 1075.5
$ws.Cells.Item(35, 10).Value = 2997.5
$ws.Cells.Item(35, 11).Value = 1075.5
$ws.Cells.Item(35, 12).Value = 2997.5
$ws.Cells.Item(35, 13).Value = -781.5
$ws.Cells.Item(35, 14).Value = -3585.5

$ws.Cells.Item(68, 8).Value = 40000
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = $null

$ws.Cells.Item(71, 8).Value = 40000
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = $null

$ws.Cells.Item(116, 8).Value = 43000
$ws.Cells.Item(116, 10).Value = 43000
$ws.Cells.Item(116, 12).Value = 43000
$ws.Cells.Item(116, 14).Value = -52178

$ws.Cells.Item(119, 8).Value = 50000
$ws.Cells.Item(119, 10).Value = 50000
$ws.Cells.Item(119, 12).Value = 50000
$ws.Cells.Item(119, 14).Value = -59676

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 3274.5
$ws.Cells.Item(103, 10).Value = 4959.4
$ws.Cells.Item(103, 12).Value = 14878.2
$ws.Cells.Item(103, 14).Value = -16636.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 276.85715
$ws.Cells.Item(13, 9).Value = 109.5
$ws.Cells.Item(13, 11).Value = 109.5
$ws.Cells.Item(13, 13).Value = 29.5

$ws.Cells.Item(80, 8).Value = 2005
$ws.Cells.Item(80, 9).Value = 2005
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 2005
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -1007
$ws.Cells.Item(80, 14).Value = $null

$ws.Cells.Item(83, 8).Value = 2005
$ws.Cells.Item(83, 9).Value = 2005
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 10025
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -5033
$ws.Cells.Item(83, 14).Value = $null

$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 13).Value = $null

$ws.Cells.Item(102, 8).Value = 3445.2222
$ws.Cells.Item(102, 9).Value = 3429.5715
$ws.Cells.Item(102, 11).Value = 3429.5715
$ws.Cells.Item(102, 13).Value = -1807.5715

$ws.Cells.Item(107, 8).Value = 675

$ws.Cells.Item(132, 8).Value = 5643.8184
$ws.Cells.Item(132, 10).Value = 3999
$ws.Cells.Item(132, 12).Value = 11997
$ws.Cells.Item(132, 14).Value = -17057

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(31, 8).Value = 17500
$ws.Cells.Item(31, 9).Value = 17500
$ws.Cells.Item(31, 11).Value = 17500
$ws.Cells.Item(31, 13).Value = -17252

$ws.Cells.Item(55, 8).Value = 1166.3334
$ws.Cells.Item(55, 9).Value = 749.6667
$ws.Cells.Item(55, 10).Value = 1583
$ws.Cells.Item(55, 11).Value = 749.6667
$ws.Cells.Item(55, 12).Value = 1583
$ws.Cells.Item(55, 13).Value = -576.6667
$ws.Cells.Item(55, 14).Value = -1929

$ws.Cells.Item(75, 8).Value = 49998
$ws.Cells.Item(75, 10).Value = 49998
$ws.Cells.Item(75, 12).Value = 49998
$ws.Cells.Item(75, 14).Value = -51870

$ws.Cells.Item(78, 8).Value = 49998
$ws.Cells.Item(78, 10).Value = 49998
$ws.Cells.Item(78, 12).Value = 149994
$ws.Cells.Item(78, 14).Value = -159354

$ws.Cells.Item(82, 8).Value = 1112.7142
$ws.Cells.Item(82, 9).Value = 964.8333
$ws.Cells.Item(82, 10).Value = 2000
$ws.Cells.Item(82, 11).Value = 964.8333
$ws.Cells.Item(82, 12).Value = 2000
$ws.Cells.Item(82, 13).Value = -603.8333
$ws.Cells.Item(82, 14).Value = -2722

$ws.Cells.Item(85, 8).Value = 1112.7142
$ws.Cells.Item(85, 9).Value = 964.8333
$ws.Cells.Item(85, 10).Value = 2000
$ws.Cells.Item(85, 11).Value = 964.8333
$ws.Cells.Item(85, 12).Value = 2000
$ws.Cells.Item(85, 13).Value = 283.1667
$ws.Cells.Item(85, 14).Value = -4496

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 13).Value = $null
